$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.035.07"
$ws.Range("E2").Value = "  +5.57%  "
$ws.Range("D3").Value = "2.279.60"
$ws.Range("E3").Value = "  +3.16%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.644"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "65.82"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +9.11%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  +7.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.104"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +16.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.65"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +18.32%  "
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").Value = "2.618.65"
$ws.Range("E14").Value = "  +3.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.95"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.834"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.86%  "
$ws.Range("D18").Value = "2.279.37"
$ws.Range("E18").Value = "  +3.06%  "
$ws.Range("D19").Value = "43.811.95"
$ws.Range("E19").Value = "  +5.07%  "
$ws.Range("D20").Value = "0.0₃0995"
$ws.Range("E20").Value = "  +10.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "74.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.72%  "
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "261.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.88%  "
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.46%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.17%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "172.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "21.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.138"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.92%  "
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("E32").Value = "  +7.44%  "
$ws.Range("E33").Value = "  +2.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0690"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.05"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("E36").Value = "  +2.72%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.87"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.07%  "
$ws.Range("B38").Value = "THORChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.87"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.46%  "
$ws.Range("E39").Value = "  +0.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0250"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.10%  "
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.04%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.41"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.11%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0980"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.76%  "
$ws.Range("E45").Value = "  +1.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "98.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.08%  "
$ws.Range("B47").Value = "Celestia"
$ws.Range("C47").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +21.66%  "
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("D49").Value = "1.474.22"
$ws.Range("E49").Value = "  +0.52%  "
$ws.Range("E50").Value = "  +6.46%  "
$ws.Range("E51").Value = "  -14.20%  "
